$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name and title text for new date (2022-04-06 -> 2022-04-07)
$ws.Name = "Through 2022-04-07"
$ws.Range("A5").Value = "April (through 04-07)"

# Update April row (row 5) values
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 13
$ws.Range("H5").Value = 17
$ws.Range("I5").Value = 24

# Update Total row (row 6) values
$ws.Range("B6").Value = 70
$ws.Range("D6").Value = 200
$ws.Range("E6").Value = 209
$ws.Range("F6").Value = 122
$ws.Range("G6").Value = 211
$ws.Range("H6").Value = 440
$ws.Range("I6").Value = 458
